# Modify turn analysis and in cluster proportion to include multiple
# phase restriction options. This inserts three new columns (B, D, E)
# around the existing B column (which becomes column C), populating
# them with the new per-phase values taken from the target data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for each row: New-B, C (previously B, except row 6 which is new),
# New-D, New-E
$data = @{
    1  = @(1400.0, 3000.0,  4000.0,  10000.0)
    2  = @(1600.0, 3500.0,  4200.0,  10000.0)
    3  = @(2300.0, 7000.0,  5000.0,  10000.0)
    4  = @(2000.0, 3500.0,  4000.0,  10000.0)
    5  = @(1800.0, 3100.0,  5500.0,  10000.0)
    6  = @(5500.0, 10000.0, 10000.0, 10000.0)
    7  = @(4000.0, 6000.0,  6500.0,  10000.0)
    8  = @(1700.0, 3000.0,  3300.0,  10000.0)
    9  = @(7000.0, 14000.0, 17000.0, 32400.0)
    10 = @(6000.0, 7500.0,  13000.0, 32400.0)
    11 = @(8900.0, 13500.0, 13500.0, 32400.0)
    12 = @(3500.0, 7700.0,  9500.0,  32400.0)
}

# Column B already carries the original cell style; copy that formatting
# into the three newly introduced columns (C, D, E) before writing values
# so every data cell keeps consistent formatting (xlPasteFormats = -4122).
$ws.Range("B1:B12").Copy() | Out-Null
$ws.Range("C1:E12").PasteSpecial(-4122) | Out-Null

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E
}
